$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ticket rows to append (rows 27-49), columns A-G.
# Column A values are prefixed with a leading apostrophe so Excel stores
# them as literal text ("2024-05-15") instead of auto-converting the
# date-like string into a date serial number.
$rows = @(
    @('''2024-05-15', '11:28:32', '-', 'Cámara no detecta skeleton', '-', '-', '-'),
    @('''2024-05-15', '11:31:03', '-', 'Detección de sealling mal puesto', '-', '-', '-'),
    @('''2024-05-15', '11:31:07', '-', 'Cámara no detecta foams', '-', '-', '-'),
    @('''2024-05-15', '11:31:09', '-', 'AOI (malla)', '-', '-', '-'),
    @('''2024-05-15', '11:35:04', 'No coge placa', '-', '-', '-', '-'),
    @('''2024-05-15', '11:35:08', 'No atornilla clips', '-', '-', '-', '-'),
    @('''2024-05-15', '11:35:10', 'Fallo cámara visión', '-', '-', '-', '-'),
    @('''2024-05-15', '11:35:13', 'Ascensor no sube', '-', '-', '-', '-'),
    @('''2024-05-15', '11:35:16', 'No atornilla tapa', '-', '-', '-', '-'),
    @('''2024-05-15', '11:35:24', 'AOI no detecta pieza', '-', '-', '-', '-'),
    @('''2024-05-15', '11:35:49', 'Fallo fijador tapa', '-', '-', '-', '-'),
    @('''2024-05-15', '11:35:51', 'Fallo fijador tapa', '-', '-', '-', '-'),
    @('''2024-05-15', '11:35:54', 'Fallo fijador tapa', '-', '-', '-', '-'),
    @('''2024-05-15', '11:35:56', 'Fallo fijador tapa', '-', '-', '-', '-'),
    @('''2024-05-15', '11:36:02', 'No coge placa', '-', '-', '-', '-'),
    @('''2024-05-15', '11:36:07', 'Fallo en paletizador', '-', '-', '-', '-'),
    @('''2024-05-15', '11:36:14', 'Fallo en paletizador', '-', '-', '-', '-'),
    @('''2024-05-15', '11:36:22', 'Fallo tolva', '-', '-', '-', '-'),
    @('''2024-05-15', '11:36:30', 'Fallo tolva', '-', '-', '-', '-'),
    @('''2024-05-15', '11:36:34', 'Fallo tolva', '-', '-', '-', '-'),
    @('''2024-05-15', '11:45:45', '-', '-', 'Tornillo atascado', '-', '-'),
    @('''2024-05-15', '11:45:47', '-', '-', 'Power atascado en prensa, cuesta sacar', '-', '-'),
    @('''2024-05-15', '11:45:51', '-', '-', 'No lee QR', '-', '-')
)

$startRow = 27
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
